# Add three new rows of region/division data to the bottom of the sheet,
# matching the rows added in the target revision (rows 166-168).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 166: Divisão=A001, Região=MATRIZ, Nome=Matriz
$ws.Range("A166").Value = "A001"
$ws.Range("C166").Value = "MATRIZ"
$ws.Range("B166").Value = "Matriz"

# Row 167: Divisão=B006, Região=MATRIZ (Nome left blank)
$ws.Range("A167").Value = "B006"
$ws.Range("C167").Value = "MATRIZ"

# Row 168: Divisão=A174, Região=MATRIZ (Nome left blank)
$ws.Range("A168").Value = "A174"
$ws.Range("C168").Value = "MATRIZ"

# Update the view's selected/active cell to match the new bottom of the sheet
$ws.Range("H162").Select()
